$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "1.00" or
# "224.77" are not auto-converted to numbers (matches original inlineStr cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "34.407.72"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").Value = "1.837.93"
$ws.Range("E3").Value = "  +3.43%  "
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "224.77"
$ws.Range("E5").Value = "  -0.37%  "
$ws.Range("E6").Value = "  +1.48%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").Value = "31.98"
$ws.Range("E8").Value = "  +0.72%  "
$ws.Range("D9").Value = "0.292"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").Value = "0.0709"
$ws.Range("E10").Value = "  +7.86%  "
$ws.Range("D11").Value = "0.0935"
$ws.Range("E11").Value = "  +0.78%  "
$ws.Range("D12").Value = "2.108.34"
$ws.Range("E12").Value = "  +3.82%  "
$ws.Range("D13").Value = "1.842.93"
$ws.Range("E13").Value = "  +3.69%  "
$ws.Range("D14").Value = "10.82"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("D15").Value = "0.647"
$ws.Range("E15").Value = "  +3.46%  "
$ws.Range("D16").Value = "34.455.79"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "4.35"
$ws.Range("E17").Value = "  +3.32%  "
$ws.Range("D18").Value = "69.68"
$ws.Range("E18").Value = "  +1.21%  "
$ws.Range("D19").Value = "251.05"
$ws.Range("E19").Value = "  -1.27%  "
$ws.Range("D20").Value = "0.0₃0795"
$ws.Range("E20").Value = "  +7.71%  "
$ws.Range("D21").Value = "11.27"
$ws.Range("E21").Value = "  +8.76%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "4.27"
$ws.Range("E23").Value = "  +1.80%  "
$ws.Range("E24").Value = "  +1.31%  "
$ws.Range("D25").Value = "161.42"
$ws.Range("E25").Value = "  +3.05%  "
$ws.Range("D26").Value = "16.64"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Value = "7.24"
$ws.Range("E27").Value = "  +3.57%  "
$ws.Range("D28").Value = "0.115"
$ws.Range("E28").Value = "  +1.68%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").Value = "0.0534"
$ws.Range("E30").Value = "  +4.64%  "
$ws.Range("D31").Value = "3.81"
$ws.Range("E31").Value = "  +1.17%  "
$ws.Range("D32").Value = "1.21"
$ws.Range("E32").Value = "  +1.34%  "
$ws.Range("D33").Value = "3.61"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").Value = "1.93"
$ws.Range("E34").Value = "  +4.23%  "
$ws.Range("D35").Value = "1.455.16"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "0.646"
$ws.Range("E36").Value = "  +3.33%  "
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "0.0193"
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("D39").Value = "0.967"
$ws.Range("E39").Value = "  +8.55%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").Value = "82.19"
$ws.Range("E40").Value = "  -0.81%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D41").Value = "2.78"
$ws.Range("E41").Value = "  -2.37%  "
$ws.Range("D42").Value = "2.36"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("E43").Value = "  +4.88%  "
$ws.Range("D44").Value = "6.11"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("D45").Value = "2.001.93"
$ws.Range("E45").Value = "  +3.62%  "
$ws.Range("E46").Value = "  +1.16%  "
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "106.78"
$ws.Range("E48").Value = "  +8.66%  "
$ws.Range("D49").Value = "12.17"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("E50").Value = "  +0.37%  "
$ws.Range("E51").Value = "  +7.88%  "

# Restore original (unset) formatting on column D so styles match the source file.
$ws.Range("D2:D51").ClearFormats()

